# Generate Report for Handoff
# Updates the "Latest Handoff Date"/"Latest Handoff Datetime" for the last
# (4th) file row in all three sheets, reflecting a new handoff generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row 5 is b1d76ec3-3e9b-4226-ab79-56cc1be0a550.md, column D
# is "Latest Handoff Date".
$wsOverview.Range("D5").Value = "2016-03-25 09:28:21"

# zh-cn sheet: row 5 is the same file; column E is "Latest Handoff Datetime".
$wsZhCn.Range("E5").Value = "2016-03-25 09:28:08"

# de-de sheet: row 5 is the same file; column E is "Latest Handoff Datetime".
$wsDeDe.Range("E5").Value = "2016-03-25 09:28:21"
